$p = $ppt.ActivePresentation

# Slide 4 ("Review of Continuations" title slide) -> retitle for the
# activity-oriented final day ("The Continuation Game").
$s = $p.Slides.Item(4)

# The title placeholder is shape 1 ("Rectangle 2"); its text is currently
# split across two runs ("Review of Continuations" + a trailing space).
# Replacing the whole TextRange collapses it back into a single run while
# preserving the first run's formatting (sz=4000, bold).
$titleShape = $s.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "The Continuation Game"
